$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "Alternant" (was F1) and "Langue" (was I1) columns are removed from the
# student sheet. "Groupe" and "Demi-groupe" slide left into F1/G1, and the
# now-empty trailing cells (old H1/I1) are cleared so the sheet shrinks to
# A1:G3 while the (now unused) column width formatting for G:H is left as-is.
$ws.Range("F1").Value = "Groupe"
$ws.Range("G1").Value = "Demi-groupe"
$ws.Range("H1:I1").ClearContents()

# Selection moves to D1 (eMail column) as the active cell.
$ws.Range("D1").Select() | Out-Null
